$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reactions")

# Insert a new "Rate units" column between "Reversible" (E) and "Flux min" (old F)
$ws.Columns("F:F").Insert() | Out-Null
$ws.Range("F1").Value = "Rate units"
$ws.Range("F2:F6").Value = "s^-1"

# Rebuild the AutoFilter so it covers the newly inserted column
$ws.AutoFilterMode = $false
$ws.Range("A1:J6").AutoFilter() | Out-Null

# The special _FilterDatabase defined names don't auto-expand with the
# inserted column, so update them explicitly to match the new range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Reactions!_FilterDatabase") {
        $n.RefersTo = "=Reactions!`$A`$1:`$J`$6"
    } elseif ($n.Name -like "Reactions!_FilterDatabase_0*") {
        $n.RefersTo = "=Reactions!`$A`$1:`$J`$1"
    }
}

# Select the newly added data and make Reactions the active sheet/tab
$ws.Range("F2:F6").Select() | Out-Null
$ws.Activate() | Out-Null
